$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Experimental (B7): blank -> "true" (kept as TEXT, not boolean) ---
# Assigning the literal string "true"/"TRUE" to a Range.Value auto-coerces
# to a native Excel boolean (t="b"). To preserve it as text (t="s", same
# cell style as neighbouring cells), stage it in a scratch cell using a
# leading apostrophe (forces text entry, stored without the apostrophe),
# then copy/paste the *value* into the target cell and clean the scratch
# cell back up so the sheet's used range is unaffected.
$scratch = $ws.Range("ZZ1")
$scratch.Value = "'true"
$scratch.Copy()
$ws.Range("B7").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

# --- Date (B8): updated timestamp (plain text, no coercion risk) ---
$ws.Range("B8").Value = "2023-02-16T14:43:10-06:00"

$excel.CutCopyMode = 0
